$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new record at row 141 (pushes old rows 141..236 down to 142..237) ---
$ws.Rows("141:141").Insert()
$ws.Cells.Item(141,1).Value = 6
$ws.Cells.Item(141,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(141,3).Value = "Metropolitana"
$ws.Cells.Item(141,4).Value = 44818
$ws.Cells.Item(141,5).Value = 13
$ws.Cells.Item(141,6).Value = 100112029
$ws.Cells.Item(141,7).Value = "Orégano"
$ws.Cells.Item(141,8).Value = "Sin especificar"
$ws.Cells.Item(141,9).Value = "Primera"
$ws.Cells.Item(141,10).Value = 46
$ws.Cells.Item(141,11).Value = 15000
$ws.Cells.Item(141,12).Value = 16000
$ws.Cells.Item(141,13).Value = 15457
$ws.Cells.Item(141,14).Value = "`$/docena de atados"
$ws.Cells.Item(141,15).Value = "Región Metropolitana"
$ws.Cells.Item(141,16).Value = 5152
$ws.Cells.Item(141,17).Value = 3
$ws.Cells.Item(141,18).Value = "Hortaliza"

# --- Insert new record at row 208 (pushes current rows 208..237 down to 209..238) ---
$ws.Rows("208:208").Insert()
$ws.Cells.Item(208,1).Value = 6
$ws.Cells.Item(208,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(208,3).Value = "Metropolitana"
$ws.Cells.Item(208,4).Value = 44816
$ws.Cells.Item(208,5).Value = 13
$ws.Cells.Item(208,6).Value = 100112029
$ws.Cells.Item(208,7).Value = "Orégano"
$ws.Cells.Item(208,8).Value = "Sin especificar"
$ws.Cells.Item(208,9).Value = "Primera"
$ws.Cells.Item(208,10).Value = 49
$ws.Cells.Item(208,11).Value = 15000
$ws.Cells.Item(208,12).Value = 16000
$ws.Cells.Item(208,13).Value = 15449
$ws.Cells.Item(208,14).Value = "`$/docena de atados"
$ws.Cells.Item(208,15).Value = "Región Metropolitana"
$ws.Cells.Item(208,16).Value = 5150
$ws.Cells.Item(208,17).Value = 3
$ws.Cells.Item(208,18).Value = "Hortaliza"
